$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the p-value used for significance (column E, row 2): 10% -> 5%
$ws.Range("E2").Value = 0.05

# Update source (K2), link (L2), notes (M2) text for the Pawlowski et al. (2019) entry
$ws.Range("K2").Value = "Pawlowski et al. (2019) Figure 2 & Text p. 19, Figure 5 of the published version"
$ws.Range("L2").Value = "https://ideas.repec.org/p/usg/econwp/201906.html"
$ws.Range("M2").Value = "The authors do not present their estimates in a table with standard errors. Instead they show the effect for each year in Figure 2. The effect is also somewhat jumpy between years.  However, the authors mention in the text what they believe to be the average effect. (i.e. 260€ for men / 0 for women when comparing high to low expenditure). Looking at the graph these appear to be significant at at about 5%. The results appear to be unchanged in the version that was published in Labor Economics in 2021."

# Increase row height to fit the longer notes text
$ws.Rows.Item(2).RowHeight = 135

# Update the selected cell in the sheet view
$ws.Range("E2").Select()
